# Expenses_details.xlsx — add new expense rows, replace the old sample
# rows, and extend the used range from A1:C4 to A1:C9.
#
# Column C holds dates stored as plain text (not real Excel dates), so
# each value is entered with a leading apostrophe (forces text / avoids
# Excel's automatic "looks like a date -> convert to date serial" coercion)
# and then ClearFormats() strips the quote-prefix cell style back off so
# the cell is left with the default style, matching the original file's
# formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-ExpenseRow([int]$Row, [string]$Category, [double]$Amount, [string]$Date) {
    $ws.Range("A$Row").Value = $Category
    $ws.Range("B$Row").Value = $Amount

    $cell = $ws.Range("C$Row")
    $cell.Value = "'" + $Date
    $cell.ClearFormats()
}

Set-ExpenseRow 2 "Books"           100  "2025-07-30"
Set-ExpenseRow 3 "Transport"       600  "2025-07-29"
Set-ExpenseRow 4 "Food"            300  "2025-07-27"
Set-ExpenseRow 5 "Transport"       120  "2025-07-25"
Set-ExpenseRow 6 "Entertainment"   450  "2025-07-20"
Set-ExpenseRow 7 "Bills"           1000 "2025-07-10"
Set-ExpenseRow 8 "IceCream"        30   "2025-07-06"
Set-ExpenseRow 9 "ElectricityBill" 500  "2025-07-06"
